$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 155 (existing rows 155-207 shift down to 157-209)
$ws.Rows("155:156").Insert()

# ---- New row 155 ----
$ws.Cells.Item(155, 1).Value2 = 9
$ws.Cells.Item(155, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(155, 3).Value2 = "Metropolitana"
$ws.Cells.Item(155, 4).Value2 = 44825
$ws.Cells.Item(155, 5).Value2 = 13
$ws.Cells.Item(155, 6).Value2 = "Fruta"
$ws.Cells.Item(155, 7).Value2 = 100101
$ws.Cells.Item(155, 8).Value2 = "Berries"
$ws.Cells.Item(155, 9).Value2 = 100101001
$ws.Cells.Item(155, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(155, 11).Value2 = "Sin especificar"
$ws.Cells.Item(155, 12).Value2 = "Primera"
$ws.Cells.Item(155, 13).Value2 = 712
$ws.Cells.Item(155, 14).Value2 = 6000
$ws.Cells.Item(155, 15).Value2 = 7000
$ws.Cells.Item(155, 16).Value2 = 6500
$ws.Cells.Item(155, 17).Value2 = "`$/bandeja 12 canastillos 125 gramos"
$ws.Cells.Item(155, 18).Value2 = "Perú"
$ws.Cells.Item(155, 19).Value2 = 4333
$ws.Cells.Item(155, 20).Value2 = 1.5

# ---- New row 156 ----
$ws.Cells.Item(156, 1).Value2 = 9
$ws.Cells.Item(156, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(156, 3).Value2 = "Metropolitana"
$ws.Cells.Item(156, 4).Value2 = 44825
$ws.Cells.Item(156, 5).Value2 = 13
$ws.Cells.Item(156, 6).Value2 = "Fruta"
$ws.Cells.Item(156, 7).Value2 = 100101
$ws.Cells.Item(156, 8).Value2 = "Berries"
$ws.Cells.Item(156, 9).Value2 = 100101001
$ws.Cells.Item(156, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(156, 11).Value2 = "Sin especificar"
$ws.Cells.Item(156, 12).Value2 = "Primera"
$ws.Cells.Item(156, 13).Value2 = 150
$ws.Cells.Item(156, 14).Value2 = 12000
$ws.Cells.Item(156, 15).Value2 = 14000
$ws.Cells.Item(156, 16).Value2 = 13000
$ws.Cells.Item(156, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(156, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(156, 19).Value2 = 6500
$ws.Cells.Item(156, 20).Value2 = 2
